# Ajout de l'activité du jour : on ajoute une nouvelle ligne (31) au journal
# de bord avec la date, le compte rendu de l'entretien avec M. Egger et le
# nombre de périodes passées.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$row = 31

# Colonne A : date du jour (2 mars 2017). On copie d'abord la mise en forme
# de la cellule de date précédente (A30) pour réutiliser le même style
# (format de date) plutôt que d'en créer un nouveau.
$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item($row, 1).Value = [DateTime]"2017-03-02"

# Colonne B : description de l'activité, avec le même style de retour à la
# ligne (wrap text) que les autres cellules de description (copié de B30).
$ws.Cells.Item($row - 1, 2).Copy() | Out-Null
$ws.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Cells.Item($row, 2).Value = "Compte rendu fait par M. Egger sur ma documentation de projet. Nous en avons discuté et j'ai commencé à améliorer les points qui sont sortis durant la disscusion"

# Colonne C : nombre de périodes travaillées
$ws.Cells.Item($row, 3).Value = "1 périodes"

$excel.CutCopyMode = 0

# Hauteur de ligne pour le texte qui s'affiche sur deux lignes, comme les
# autres lignes de description (18, 19, 24, 30, ...)
$ws.Rows.Item($row).RowHeight = 30

# La sélection active se déplace sur la cellule suivante, comme reflété
# dans la vue de la feuille après l'ajout de la nouvelle ligne.
$ws.Range("C32").Select()
